$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the two trailing rows (testuser3 / testuser4) ---
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# --- Drop the existing hyperlinks so we can rebuild just the two we keep ---
$ws.Hyperlinks.Delete()

# --- Row 2: stevewauhghg / user_davisjames@123 / male / inactive ---
$ws.Range("A2").Value = "stevewauhghg"
$ws.Range("B2").Value = "user_davisjames@123"
$ws.Range("C2").Value = "male"
$ws.Range("D2").Value = "inactive"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:user_davisjames@123")
$ws.Range("B2").Style = "Hyperlink"

# --- Row 3: janesgerde / user_scottstevy@456 / female / active ---
$ws.Range("A3").Value = "janesgerde"
$ws.Range("B3").Value = "user_scottstevy@456"
$ws.Range("C3").Value = "female"
$ws.Range("D3").Value = "active"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:user_scottstevy@456")
$ws.Range("B3").Style = "Hyperlink"

# --- Column widths (closest values this engine's 1/MDW rounding can reach) ---
$ws.Columns.Item(1).ColumnWidth = 17.833333333333332
$ws.Columns.Item(2).ColumnWidth = 35.0
$ws.Columns.Item(3).ColumnWidth = 17.666666666666668
$ws.Columns.Item(4).ColumnWidth = 17.5

# --- Selection moves to D3 ---
$ws.Range("D3").Select() | Out-Null
